$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ly9"
$ws.Cells.Item(2,3).Value = "Ly9"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1213926666666667
$ws.Cells.Item(2,8).Value = 0.364178
$ws.Cells.Item(2,9).Value = 0.001287165941041644
$ws.Cells.Item(2,10).Value = 0.001287165941041644
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1213926666666667
$ws.Cells.Item(2,14).Value = 0.364178
$ws.Cells.Item(2,15).Value = 0.001287165941041644
$ws.Cells.Item(2,16).Value = 0.001287165941041644
$ws.Cells.Item(2,17).Value = 0.01473617952044444
$ws.Cells.Item(2,18).Value = 0.132625615684
$ws.Cells.Item(2,19).Value = 0.000001656796159777622
$ws.Cells.Item(2,20).Value = 0.000001656796159777621

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ly9"
$ws.Cells.Item(3,3).Value = "Ly9"
$ws.Cells.Item(3,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1213926666666667
$ws.Cells.Item(3,8).Value = 0.364178
$ws.Cells.Item(3,9).Value = 0.001287165941041644
$ws.Cells.Item(3,10).Value = 0.001287165941041644
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 49.81377800000001
$ws.Cells.Item(3,14).Value = 149.441334
$ws.Cells.Item(3,15).Value = 0.5281916955681801
$ws.Cells.Item(3,16).Value = 0.52819169556818
$ws.Cells.Item(3,17).Value = 6.047027348161333
$ws.Cells.Item(3,18).Value = 54.423246133452
$ws.Cells.Item(3,19).Value = 0.0006798703608763984
$ws.Cells.Item(3,20).Value = 0.000679870360876398

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ly9"
$ws.Cells.Item(4,3).Value = "Ly9"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1213926666666667
$ws.Cells.Item(4,8).Value = 0.364178
$ws.Cells.Item(4,9).Value = 0.001287165941041644
$ws.Cells.Item(4,10).Value = 0.001287165941041644
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 44.37486566666666
$ws.Cells.Item(4,14).Value = 133.124597
$ws.Cells.Item(4,15).Value = 0.4705211384907783
$ws.Cells.Item(4,16).Value = 0.4705211384907783
$ws.Cells.Item(4,17).Value = 5.386783276251777
$ws.Cells.Item(4,18).Value = 48.481049486266
$ws.Cells.Item(4,19).Value = 0.0006056387840054685
$ws.Cells.Item(4,20).Value = 0.0006056387840054684

# Row 5
$ws.Cells.Item(5,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(5,2).Value = "Ly9"
$ws.Cells.Item(5,3).Value = "Ly9"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 49.81377800000001
$ws.Cells.Item(5,8).Value = 149.441334
$ws.Cells.Item(5,9).Value = 0.5281916955681801
$ws.Cells.Item(5,10).Value = 0.52819169556818
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.1213926666666667
$ws.Cells.Item(5,14).Value = 0.364178
$ws.Cells.Item(5,15).Value = 0.001287165941041644
$ws.Cells.Item(5,16).Value = 0.001287165941041644
$ws.Cells.Item(5,17).Value = 6.047027348161333
$ws.Cells.Item(5,18).Value = 54.423246133452
$ws.Cells.Item(5,19).Value = 0.0006798703608763984
$ws.Cells.Item(5,20).Value = 0.000679870360876398

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Ly9"
$ws.Cells.Item(6,3).Value = "Ly9"
$ws.Cells.Item(6,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 49.81377800000001
$ws.Cells.Item(6,8).Value = 149.441334
$ws.Cells.Item(6,9).Value = 0.5281916955681801
$ws.Cells.Item(6,10).Value = 0.52819169556818
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 49.81377800000001
$ws.Cells.Item(6,14).Value = 149.441334
$ws.Cells.Item(6,15).Value = 0.5281916955681801
$ws.Cells.Item(6,16).Value = 0.52819169556818
$ws.Cells.Item(6,17).Value = 2481.412478633285
$ws.Cells.Item(6,18).Value = 22332.71230769956
$ws.Cells.Item(6,19).Value = 0.2789864672671891
$ws.Cells.Item(6,20).Value = 0.278986467267189

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Ly9"
$ws.Cells.Item(7,3).Value = "Ly9"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 49.81377800000001
$ws.Cells.Item(7,8).Value = 149.441334
$ws.Cells.Item(7,9).Value = 0.5281916955681801
$ws.Cells.Item(7,10).Value = 0.52819169556818
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 44.37486566666666
$ws.Cells.Item(7,14).Value = 133.124597
$ws.Cells.Item(7,15).Value = 0.4705211384907783
$ws.Cells.Item(7,16).Value = 0.4705211384907783
$ws.Cells.Item(7,17).Value = 2210.479707099155
$ws.Cells.Item(7,18).Value = 19894.3173638924
$ws.Cells.Item(7,19).Value = 0.2485253579401147
$ws.Cells.Item(7,20).Value = 0.2485253579401146

# Row 8
$ws.Cells.Item(8,1).Value = "Resolving-Mac"
$ws.Cells.Item(8,2).Value = "Ly9"
$ws.Cells.Item(8,3).Value = "Ly9"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 44.37486566666666
$ws.Cells.Item(8,8).Value = 133.124597
$ws.Cells.Item(8,9).Value = 0.4705211384907783
$ws.Cells.Item(8,10).Value = 0.4705211384907783
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.1213926666666667
$ws.Cells.Item(8,14).Value = 0.364178
$ws.Cells.Item(8,15).Value = 0.001287165941041644
$ws.Cells.Item(8,16).Value = 0.001287165941041644
$ws.Cells.Item(8,17).Value = 5.386783276251777
$ws.Cells.Item(8,18).Value = 48.481049486266
$ws.Cells.Item(8,19).Value = 0.0006056387840054685
$ws.Cells.Item(8,20).Value = 0.0006056387840054684

# Row 9
$ws.Cells.Item(9,1).Value = "Resolving-Mac"
$ws.Cells.Item(9,2).Value = "Ly9"
$ws.Cells.Item(9,3).Value = "Ly9"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 44.37486566666666
$ws.Cells.Item(9,8).Value = 133.124597
$ws.Cells.Item(9,9).Value = 0.4705211384907783
$ws.Cells.Item(9,10).Value = 0.4705211384907783
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 49.81377800000001
$ws.Cells.Item(9,14).Value = 149.441334
$ws.Cells.Item(9,15).Value = 0.5281916955681801
$ws.Cells.Item(9,16).Value = 0.52819169556818
$ws.Cells.Item(9,17).Value = 2210.479707099155
$ws.Cells.Item(9,18).Value = 19894.3173638924
$ws.Cells.Item(9,19).Value = 0.2485253579401147
$ws.Cells.Item(9,20).Value = 0.2485253579401146

# Row 10
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Ly9"
$ws.Cells.Item(10,3).Value = "Ly9"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 44.37486566666666
$ws.Cells.Item(10,8).Value = 133.124597
$ws.Cells.Item(10,9).Value = 0.4705211384907783
$ws.Cells.Item(10,10).Value = 0.4705211384907783
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 44.37486566666666
$ws.Cells.Item(10,14).Value = 133.124597
$ws.Cells.Item(10,15).Value = 0.4705211384907783
$ws.Cells.Item(10,16).Value = 0.4705211384907783
$ws.Cells.Item(10,17).Value = 1969.128702934712
$ws.Cells.Item(10,18).Value = 17722.15832641241
$ws.Cells.Item(10,19).Value = 0.2213901417666582
$ws.Cells.Item(10,20).Value = 0.2213901417666581
